$wb = $excel.ActiveWorkbook

# --- Update the summary text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.25 = 12370.13 pesos`n✅ 12370.13 pesos = 3.22 = 969.53 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 308
$wsTasas.Range("O10").Value = 3810
$wsTasas.Range("N12").Value = 3838.5
$wsTasas.Range("O12").Value = 300.85
